$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($range, $values)
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($range).Value = $arr
}

# Update existing rows 2-5 with new schedule data (columns B:I)
Set-RowValues "B2:I2" @(8, 2, 4, 4, -4, 2, 23, 5)
Set-RowValues "B3:I3" @(7, 2, 2, 3, -5, 1, 12, 5)
Set-RowValues "B4:I4" @(6, 3, 5, 8, -1, 5, 56, 5)
Set-RowValues "B5:I5" @(9, 4, 7, 8, -2, 4, 45, 5)

# Add new row 6 with a new trial
$ws.Range("A6").Value = 5
Set-RowValues "B6:I6" @(7, 0, 4, 3, -3, 3, 34, 5)
$ws.Range("J6").Value = "train_dim2_1"

$ws.Range("I1").Select()
